$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 184.0626906666667
$ws.Cells.Item(2, 8).Value = 552.188072
$ws.Cells.Item(2, 9).Value = 0.6510505751503485
$ws.Cells.Item(2, 10).Value = 0.6510505751503486
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 505.8908573333333
$ws.Cells.Item(2, 14).Value = 1517.672572
$ws.Cells.Item(2, 15).Value = 0.7018211771568338
$ws.Cells.Item(2, 16).Value = 0.7018211771568337
$ws.Cells.Item(2, 17).Value = 93115.63238444013
$ws.Cells.Item(2, 18).Value = 838040.6914599612
$ws.Cells.Item(2, 19).Value = 0.4569210810406513
$ws.Cells.Item(2, 20).Value = 0.4569210810406513

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 184.0626906666667
$ws.Cells.Item(3, 8).Value = 552.188072
$ws.Cells.Item(3, 9).Value = 0.6510505751503485
$ws.Cells.Item(3, 10).Value = 0.6510505751503486
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 88.00803400000001
$ws.Cells.Item(3, 14).Value = 264.024102
$ws.Cells.Item(3, 15).Value = 0.1220933352041998
$ws.Cells.Item(3, 16).Value = 0.1220933352041997
$ws.Cells.Item(3, 17).Value = 16198.99553832348
$ws.Cells.Item(3, 18).Value = 145790.9598449114
$ws.Cells.Item(3, 19).Value = 0.07948893610671853
$ws.Cells.Item(3, 20).Value = 0.07948893610671855

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 184.0626906666667
$ws.Cells.Item(4, 8).Value = 552.188072
$ws.Cells.Item(4, 9).Value = 0.6510505751503485
$ws.Cells.Item(4, 10).Value = 0.6510505751503486
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 126.926974
$ws.Cells.Item(4, 14).Value = 380.780922
$ws.Cells.Item(4, 15).Value = 0.1760854876389666
$ws.Cells.Item(4, 16).Value = 0.1760854876389665
$ws.Cells.Item(4, 17).Value = 23362.52035261804
$ws.Cells.Item(4, 18).Value = 210262.6831735624
$ws.Cells.Item(4, 19).Value = 0.1146405580029788
$ws.Cells.Item(4, 20).Value = 0.1146405580029787

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 57.4434
$ws.Cells.Item(5, 8).Value = 172.3302
$ws.Cells.Item(5, 9).Value = 0.2031838091312023
$ws.Cells.Item(5, 10).Value = 0.2031838091312023
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 505.8908573333333
$ws.Cells.Item(5, 14).Value = 1517.672572
$ws.Cells.Item(5, 15).Value = 0.7018211771568338
$ws.Cells.Item(5, 16).Value = 0.7018211771568337
$ws.Cells.Item(5, 17).Value = 29060.0908741416
$ws.Cells.Item(5, 18).Value = 261540.8178672744
$ws.Cells.Item(5, 19).Value = 0.1425987001036698
$ws.Cells.Item(5, 20).Value = 0.1425987001036698

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 57.4434
$ws.Cells.Item(6, 8).Value = 172.3302
$ws.Cells.Item(6, 9).Value = 0.2031838091312023
$ws.Cells.Item(6, 10).Value = 0.2031838091312023
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 88.00803400000001
$ws.Cells.Item(6, 14).Value = 264.024102
$ws.Cells.Item(6, 15).Value = 0.1220933352041998
$ws.Cells.Item(6, 16).Value = 0.1220933352041997
$ws.Cells.Item(6, 17).Value = 5055.480700275601
$ws.Cells.Item(6, 18).Value = 45499.3263024804
$ws.Cells.Item(6, 19).Value = 0.02480738891632202
$ws.Cells.Item(6, 20).Value = 0.02480738891632202

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 57.4434
$ws.Cells.Item(7, 8).Value = 172.3302
$ws.Cells.Item(7, 9).Value = 0.2031838091312023
$ws.Cells.Item(7, 10).Value = 0.2031838091312023
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 126.926974
$ws.Cells.Item(7, 14).Value = 380.780922
$ws.Cells.Item(7, 15).Value = 0.1760854876389666
$ws.Cells.Item(7, 16).Value = 0.1760854876389665
$ws.Cells.Item(7, 17).Value = 7291.116938271601
$ws.Cells.Item(7, 18).Value = 65620.0524444444
$ws.Cells.Item(7, 19).Value = 0.03577772011121046
$ws.Cells.Item(7, 20).Value = 0.03577772011121045

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 41.21033366666666
$ws.Cells.Item(8, 8).Value = 123.631001
$ws.Cells.Item(8, 9).Value = 0.1457656157184491
$ws.Cells.Item(8, 10).Value = 0.1457656157184491
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 505.8908573333333
$ws.Cells.Item(8, 14).Value = 1517.672572
$ws.Cells.Item(8, 15).Value = 0.7018211771568338
$ws.Cells.Item(8, 16).Value = 0.7018211771568337
$ws.Cells.Item(8, 17).Value = 20847.93102962273
$ws.Cells.Item(8, 18).Value = 187631.3792666045
$ws.Cells.Item(8, 19).Value = 0.1023013960125126
$ws.Cells.Item(8, 20).Value = 0.1023013960125126

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 41.21033366666666
$ws.Cells.Item(9, 8).Value = 123.631001
$ws.Cells.Item(9, 9).Value = 0.1457656157184491
$ws.Cells.Item(9, 10).Value = 0.1457656157184491
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 88.00803400000001
$ws.Cells.Item(9, 14).Value = 264.024102
$ws.Cells.Item(9, 15).Value = 0.1220933352041998
$ws.Cells.Item(9, 16).Value = 0.1220933352041997
$ws.Cells.Item(9, 17).Value = 3626.840446487345
$ws.Cells.Item(9, 18).Value = 32641.5640183861
$ws.Cells.Item(9, 19).Value = 0.01779701018115918
$ws.Cells.Item(9, 20).Value = 0.01779701018115917

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 41.21033366666666
$ws.Cells.Item(10, 8).Value = 123.631001
$ws.Cells.Item(10, 9).Value = 0.1457656157184491
$ws.Cells.Item(10, 10).Value = 0.1457656157184491
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 126.926974
$ws.Cells.Item(10, 14).Value = 380.780922
$ws.Cells.Item(10, 15).Value = 0.1760854876389666
$ws.Cells.Item(10, 16).Value = 0.1760854876389665
$ws.Cells.Item(10, 17).Value = 5230.702949840325
$ws.Cells.Item(10, 18).Value = 47076.32654856292
$ws.Cells.Item(10, 19).Value = 0.02566720952477732
$ws.Cells.Item(10, 20).Value = 0.02566720952477732
